$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format before writing, so numeric-looking
# strings (e.g. "516.42") and percentage strings stay text like the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "57.566.16"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.084.38"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "516.42"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "142.90"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "7.29"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "3.614.58"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "25.75"
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "57.660.05"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.086.03"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "6.15"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "337.41"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "65.59"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "0.0₃0927"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "20.92"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").Value = "153.92"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "27.86"
$ws.Range("E34").Value = "  +8.31%  "
$ws.Range("D35").Value = "4.52"
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "0.0689"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").Value = "3.122.81"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "36.84"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "0.671"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "2.285.16"
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "0.948"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "20.29"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "5.89"
$ws.Range("E49").Value = "  -5.35%  "
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E51").Value = "  +1.17%  "

# Restore the default (unstyled) look for columns D and E, matching the source.
$ws.Range("D2:E51").Style = "Normal"
